# Update generated 广州-漫展信息.xlsx data (commit: "Update gh-pages to
# output generated at 456a3b4") - refreshes the "想去人数" (F) counters and a
# couple of "最低票价" (G) sold-out labels across all four sheets.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item(1)   # 展览
$ws2 = $wb.Worksheets.Item(2)   # 演出
$ws3 = $wb.Worksheets.Item(3)   # 本地生活
$ws4 = $wb.Worksheets.Item(4)   # 全部类型

# ---- 展览 (sheet 1) ----
$ws1.Range("F2").Value = 2367
$ws1.Range("F3").Value = 553
$ws1.Range("F4").Value = 211
$ws1.Range("F5").Value = 361
$ws1.Range("F6").Value = 361
$ws1.Range("F7").Value = 593
$ws1.Range("F9").Value = 804
$ws1.Range("F11").Value = 833
$ws1.Range("F14").Value = 403
$ws1.Range("F15").Value = 21
$ws1.Range("F16").Value = 1033
$ws1.Range("F17").Value = 21526
$ws1.Range("G17").Value = "已售罄"
$ws1.Range("F18").Value = 893
$ws1.Range("F19").Value = 83
$ws1.Range("F20").Value = 275
$ws1.Range("F21").Value = 308
$ws1.Range("F22").Value = 180
$ws1.Range("F23").Value = 171
$ws1.Range("F25").Value = 20
$ws1.Range("F26").Value = 253
$ws1.Range("F28").Value = 364
$ws1.Range("F29").Value = 162

# ---- 演出 (sheet 2) ----
$ws2.Range("F6").Value = 209
$ws2.Range("F7").Value = 230
$ws2.Range("F8").Value = 3450
$ws2.Range("F16").Value = 3927

# ---- 本地生活 (sheet 3) ----
$ws3.Range("F3").Value = 120
$ws3.Range("F4").Value = 637
$ws3.Range("F5").Value = 216

# ---- 全部类型 (sheet 4) ----
$ws4.Range("F3").Value = 120
$ws4.Range("F5").Value = 2367
$ws4.Range("F6").Value = 637
$ws4.Range("F7").Value = 553
$ws4.Range("F8").Value = 211
$ws4.Range("F9").Value = 361
$ws4.Range("F10").Value = 361
$ws4.Range("F11").Value = 593
$ws4.Range("F16").Value = 209
$ws4.Range("F17").Value = 216
$ws4.Range("F18").Value = 804
$ws4.Range("F20").Value = 833
$ws4.Range("F23").Value = 403
$ws4.Range("F24").Value = 21
$ws4.Range("F25").Value = 1033
$ws4.Range("F26").Value = 21526
$ws4.Range("G26").Value = 0
$ws4.Range("F27").Value = 230
$ws4.Range("F28").Value = 3450
$ws4.Range("F32").Value = 893
$ws4.Range("F33").Value = 83
$ws4.Range("F34").Value = 275
$ws4.Range("F37").Value = 308
$ws4.Range("F38").Value = 180
$ws4.Range("F39").Value = 171
$ws4.Range("F41").Value = 20
$ws4.Range("F44").Value = 253
$ws4.Range("F46").Value = 364
$ws4.Range("F47").Value = 162
$ws4.Range("F48").Value = 3927
